$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.916968030515761
$ws.Range("C2").Value = 0.4206471542703127
$ws.Range("D2").Value = 0.05283817609281094
$ws.Range("F2").Value = 4.569849501827576
$ws.Range("G2").Value = 0.002543688357559137
$ws.Range("J2").Value = 0.3550245360243167
$ws.Range("B3").Value = 1.815928658003202
$ws.Range("C3").Value = 0.3932419509856686
$ws.Range("D3").Value = 0.05426609883985556
$ws.Range("F3").Value = 4.391096460491724
$ws.Range("G3").Value = 0.002550902183112819
$ws.Range("J3").Value = 0.3445332050459058
$ws.Range("B4").Value = 1.755608591331111
$ws.Range("C4").Value = 0.3768376215979004
$ws.Range("D4").Value = 0.05520331806364531
$ws.Range("F4").Value = 4.28321579520653
$ws.Range("G4").Value = 0.00255555492809739
$ws.Range("J4").Value = 0.3383344652334586
$ws.Range("B5").Value = 1.731456115233584
$ws.Range("C5").Value = 0.3702577913089726
$ws.Range("D5").Value = 0.0556004025904322
$ws.Range("F5").Value = 4.239714173020019
$ws.Range("G5").Value = 0.002557507374214049
$ws.Range("J5").Value = 0.3358689212163597
$ws.Range("B6").Value = 1.727471385856688
$ws.Range("C6").Value = 0.3691715270662996
$ws.Range("D6").Value = 0.05566725251143012
$ws.Range("F6").Value = 4.232518271460407
$ws.Range("G6").Value = 0.002557834990466431
$ws.Range("J6").Value = 0.3354631556382941
$ws.Range("B7").Value = 1.755281132645678
$ws.Range("C7").Value = 0.3767484598433839
$ws.Range("D7").Value = 0.05520861195966731
$ws.Range("F7").Value = 4.282627266132636
$ws.Range("G7").Value = 0.002555581030683523
$ws.Range("J7").Value = 0.3383009699028321
$ws.Range("B8").Value = 1.881770684817354
$ws.Range("C8").Value = 0.4111093969792705
$ws.Range("D8").Value = 0.0533179458913402
$ws.Range("F8").Value = 4.507819877814995
$ws.Range("G8").Value = 0.002546129452874193
$ws.Range("J8").Value = 0.3513563220859055
$ws.Range("B9").Value = 2.143638978230172
$ws.Range("C9").Value = 0.4819033789033824
$ws.Range("D9").Value = 0.05009223149239261
$ws.Range("F9").Value = 4.964804713883211
$ws.Range("G9").Value = 0.002529357084692766
$ws.Range("J9").Value = 0.378915311869747
$ws.Range("B10").Value = 2.34474894872659
$ws.Range("C10").Value = 0.5360861669600467
$ws.Range("D10").Value = 0.04801901271832598
$ws.Range("F10").Value = 5.310674695877594
$ws.Range("G10").Value = 0.002518093708430713
$ws.Range("J10").Value = 0.4004003587640454
$ws.Range("B11").Value = 2.43819769760853
$ws.Range("C11").Value = 0.5612272148414377
$ws.Range("D11").Value = 0.04714095936677865
$ws.Range("F11").Value = 5.470383752351438
$ws.Range("G11").Value = 0.002513196497235921
$ws.Range("J11").Value = 0.4104531245575771
$ws.Range("B12").Value = 2.473871669326115
$ws.Range("C12").Value = 0.5708199888647414
$ws.Range("D12").Value = 0.04681788111540541
$ws.Range("F12").Value = 5.531215428690473
$ws.Range("G12").Value = 0.002511374378877525
$ws.Range("J12").Value = 0.4143007661348435
$ws.Range("B13").Value = 2.46617580200774
$ws.Range("C13").Value = 0.5687507685775586
$ws.Range("D13").Value = 0.04688704170395397
$ws.Range("F13").Value = 5.518098325323649
$ws.Range("G13").Value = 0.00251176536974335
$ws.Range("J13").Value = 0.413470277633607
$ws.Range("B14").Value = 2.441126838126479
$ws.Range("C14").Value = 0.5620149582173894
$ws.Range("D14").Value = 0.04711419031964681
$ws.Range("F14").Value = 5.47538125822598
$ws.Range("G14").Value = 0.002513045943438502
$ws.Range("J14").Value = 0.4107688487848691
$ws.Range("B15").Value = 2.425821137372509
$ws.Range("C15").Value = 0.5578985547312527
$ws.Range("D15").Value = 0.04725455434801873
$ws.Range("F15").Value = 5.449262205104333
$ws.Range("G15").Value = 0.002513834537856943
$ws.Range("J15").Value = 0.4091194913949607
$ws.Range("B16").Value = 2.338681748977706
$ws.Range("C16").Value = 0.5344531891155384
$ws.Range("D16").Value = 0.04807771035276609
$ws.Range("F16").Value = 5.30028609738838
$ws.Range("G16").Value = 0.002518418292045193
$ws.Range("J16").Value = 0.3997490662456471
$ws.Range("B17").Value = 2.285730581308599
$ws.Range("C17").Value = 0.5201975502317282
$ws.Range("D17").Value = 0.04859940005987795
$ws.Range("F17").Value = 5.209510096509803
$ws.Range("G17").Value = 0.002521288140041926
$ws.Range("J17").Value = 0.3940726223513025
$ws.Range("B18").Value = 2.255458928533358
$ws.Range("C18").Value = 0.5120443965498112
$ws.Range("D18").Value = 0.04890558226980524
$ws.Range("F18").Value = 5.157520564986726
$ws.Range("G18").Value = 0.002522960139102609
$ws.Range("J18").Value = 0.3908339076154732
$ws.Range("B19").Value = 2.245241007009213
$ws.Range("C19").Value = 0.5092917879033507
$ws.Range("D19").Value = 0.04901029948857882
$ws.Range("F19").Value = 5.139955631132636
$ws.Range("G19").Value = 0.002523529920784051
$ws.Range("J19").Value = 0.3897418164565636
$ws.Range("B20").Value = 2.291348199157255
$ws.Range("C20").Value = 0.5217102833634044
$ws.Range("D20").Value = 0.04854323142926376
$ws.Range("F20").Value = 5.219150247943958
$ws.Range("G20").Value = 0.002520980433004594
$ws.Range("J20").Value = 0.3946741699277823
$ws.Range("B21").Value = 2.448476499734966
$ws.Range("C21").Value = 0.5639914516879116
$ws.Range("D21").Value = 0.04704721503480513
$ws.Range("F21").Value = 5.487918607558868
$ws.Range("G21").Value = 0.00251266893182599
$ws.Range("J21").Value = 0.411561208515181
$ws.Range("B22").Value = 2.552843839688876
$ws.Range("C22").Value = 0.5920474054196916
$ws.Range("D22").Value = 0.04612442978861608
$ws.Range("F22").Value = 5.665639604176533
$ws.Range("G22").Value = 0.002507425334750576
$ws.Range("J22").Value = 0.4228364870939458
$ws.Range("B23").Value = 2.496986174824087
$ws.Range("C23").Value = 0.5770342205564134
$ws.Range("D23").Value = 0.04661188760558943
$ws.Range("F23").Value = 5.570593440261121
$ws.Range("G23").Value = 0.002510206773803425
$ws.Range("J23").Value = 0.4167965745865416
$ws.Range("B24").Value = 2.288807942855101
$ws.Range("C24").Value = 0.5210262441294731
$ws.Range("D24").Value = 0.04856860578636812
$ws.Range("F24").Value = 5.214791317193686
$ws.Range("G24").Value = 0.002521119478630341
$ws.Range("J24").Value = 0.394402133238188
$ws.Range("B25").Value = 2.071287536128068
$ws.Range("C25").Value = 0.4623772575969838
$ws.Range("D25").Value = 0.05091298837378666
$ws.Range("F25").Value = 4.839455023187185
$ws.Range("G25").Value = 0.002533707351248196
$ws.Range("J25").Value = 0.3712458148970654
